# Applies the crypto price/volume refresh described by the commit diff.
# Cell values come from a scraped data refresh; D-column numeric-looking
# strings ("1.70", "0.999", ...) must stay TEXT (they are prices rendered as
# strings in the source sheet), so we force them with a leading apostrophe
# and then restore the default (unstyled) look by copying style from an
# untouched text cell (B2), since the apostrophe text-entry mints a new
# 'Text' number-format style on the cell otherwise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$styleDonor = $ws.Range("B2")

$ws.Range("D2").Value = '42.082.45'
$ws.Range("E2").Value = '  -2.03%  '
$ws.Range("D3").Value = '2.262.98'
$ws.Range("E3").Value = '  -3.32%  '
$ws.Range("D5").Value = '''299.54'
$ws.Range("D5").Style = $styleDonor.Style
$ws.Range("E5").Value = '  -2.37%  '
$ws.Range("D6").Value = '''94.22'
$ws.Range("D6").Style = $styleDonor.Style
$ws.Range("E6").Value = '  -6.30%  '
$ws.Range("D7").Value = '''0.496'
$ws.Range("D7").Style = $styleDonor.Style
$ws.Range("E7").Value = '  -3.18%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -3.95%  '
$ws.Range("D10").Value = '''33.02'
$ws.Range("D10").Style = $styleDonor.Style
$ws.Range("E10").Value = '  -5.78%  '
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").Value = '''47.98'
$ws.Range("D12").Style = $styleDonor.Style
$ws.Range("E12").Value = '  -8.01%  '
$ws.Range("D13").Value = '''0.113'
$ws.Range("D13").Style = $styleDonor.Style
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").Value = '''6.67'
$ws.Range("D14").Style = $styleDonor.Style
$ws.Range("E14").Value = '  -2.05%  '
$ws.Range("D15").Value = '2.614.84'
$ws.Range("E15").Value = '  -3.38%  '
$ws.Range("D16").Value = '''15.40'
$ws.Range("D16").Style = $styleDonor.Style
$ws.Range("E16").Value = '  -3.83%  '
$ws.Range("D17").Value = '2.263.27'
$ws.Range("E17").Value = '  -4.10%  '
$ws.Range("D18").Value = '''0.773'
$ws.Range("D18").Style = $styleDonor.Style
$ws.Range("E18").Value = '  -4.53%  '
$ws.Range("D19").Value = '42.073.31'
$ws.Range("E19").Value = '  -1.88%  '
$ws.Range("D20").Value = '0.0₃0891'
$ws.Range("E20").Value = '  -2.36%  '
$ws.Range("D21").Value = '''6.00'
$ws.Range("D21").Style = $styleDonor.Style
$ws.Range("E21").Value = '  -3.50%  '
$ws.Range("D22").Value = '''11.35'
$ws.Range("D22").Style = $styleDonor.Style
$ws.Range("E22").Value = '  -3.27%  '
$ws.Range("D23").Value = '''66.65'
$ws.Range("D23").Style = $styleDonor.Style
$ws.Range("E23").Value = '  -1.93%  '
$ws.Range("D24").Value = '''233.35'
$ws.Range("D24").Style = $styleDonor.Style
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("E25").Value = '  -5.18%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").Value = '''2.46'
$ws.Range("D27").Style = $styleDonor.Style
$ws.Range("E27").Value = '  -4.30%  '
$ws.Range("D28").Value = '''23.77'
$ws.Range("D28").Style = $styleDonor.Style
$ws.Range("E28").Value = '  -7.20%  '
$ws.Range("D29").Value = '''2.26'
$ws.Range("D29").Style = $styleDonor.Style
$ws.Range("E29").Value = '  -2.71%  '
$ws.Range("D30").Value = '''167.16'
$ws.Range("D30").Style = $styleDonor.Style
$ws.Range("E30").Value = '  +3.40%  '
$ws.Range("D31").Value = '''33.76'
$ws.Range("D31").Style = $styleDonor.Style
$ws.Range("E31").Value = '  -3.96%  '
$ws.Range("E32").Value = '  -3.33%  '
$ws.Range("D33").Value = '''0.999'
$ws.Range("D33").Style = $styleDonor.Style
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '''4.93'
$ws.Range("D34").Style = $styleDonor.Style
$ws.Range("E34").Value = '  -4.05%  '
$ws.Range("E35").Value = '  -3.89%  '
$ws.Range("E36").Value = '  -5.63%  '
$ws.Range("D37").Value = '''0.0693'
$ws.Range("D37").Style = $styleDonor.Style
$ws.Range("E37").Value = '  -4.94%  '
$ws.Range("E38").Value = '  -7.29%  '
$ws.Range("E39").Value = '  -6.07%  '
$ws.Range("D40").Value = '''0.0993'
$ws.Range("D40").Style = $styleDonor.Style
$ws.Range("E40").Value = '  -3.53%  '
$ws.Range("E41").Value = '  -3.62%  '
$ws.Range("D42").Value = '''1.70'
$ws.Range("D42").Style = $styleDonor.Style
$ws.Range("E42").Value = '  -8.77%  '
$ws.Range("E43").Value = '  -1.51%  '
$ws.Range("D44").Value = '1.956.19'
$ws.Range("E44").Value = '  -2.95%  '
$ws.Range("E45").Value = '  -2.54%  '
$ws.Range("D46").Value = '''17.44'
$ws.Range("D46").Style = $styleDonor.Style
$ws.Range("E46").Value = '  -6.93%  '
$ws.Range("E47").Value = '  -7.44%  '
$ws.Range("E48").Value = '  -4.92%  '
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("D50").Value = '2.488.19'
$ws.Range("D51").Value = '''51.76'
$ws.Range("D51").Style = $styleDonor.Style
$ws.Range("E51").Value = '  -7.28%  '
